$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update RF column (I) for rows 42 through 72 with the new recomputed value
$ws.Range("I42:I72").Value = 14.1090625
